# The deck's two embedded theme parts (ppt/theme/theme1.xml - used by the
# slide master - and ppt/theme/theme2.xml - used by the notes master) had
# their contents swapped: theme1 used to carry the "Integral" colour
# scheme and theme2 the stock "Office Theme" colour scheme; after the
# edit theme1 carries the "Office Theme" colours and theme2 carries the
# "Integral" colours.
#
# The font scheme and format scheme (fills/lines/effects) are byte-for-byte
# identical between the two themes, so the only thing that actually needs
# to change is the 12-slot colour scheme (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink). We rewrite those 12 slots on the presentation's theme
# through the ThemeColorScheme object so the master's theme part ends up
# holding the "Office Theme" palette.

function Convert-HexToOleColor {
    param([string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette ("Office Theme"), in clrScheme slot order.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$p = $ppt.ActivePresentation
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = Convert-HexToOleColor $officeThemeColors[$i - 1]
}
